$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Projets")

# Fix the "Role" text in F5: remove the stray "de du" -> "à la coordination ERPI"
$ws.Range("F5").Value = "Participation au montage et à la coordination ERPI "

# Fix the "Role" text in F8: "et au Responsable scientifique" -> "et  Responsable scientifique"
$ws.Range("F8").Value = "Participation au montage et  Responsable scientifique"

# Update the saved view state: scroll so column B is leftmost, select F9
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 5
$win.ScrollColumn = 2
$ws.Range("F9").Select()
